$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Tiles" to "FuelTanks"
$ws.Name = "FuelTanks"

# Update cell values: replace formula in A2 with a plain value, and update A1/B1/B2
$ws.Range("A1").Value = 10
$ws.Range("B1").Value = 2
$ws.Range("A2").Value = 20
$ws.Range("B2").Value = 2

# Update selection to A3
$ws.Range("A3").Select()
